# Apply Al-Hada validation edits:
# 1. Rename sheet "تقرير الأنظمة" -> "تقرير أنظمة البناء والتنظيم"
# 2. Assumptions!D3 "Fund Name " -> "Fund Name" (strip trailing space)
# 3. Assumptions!E3 "صندوق استثمار عقاري" -> "صندوق استثمار عقاري خاص"
# 4. Assumptions!D4 "Fund Type " -> "Fund Type" (strip trailing space)

$wb = $excel.ActiveWorkbook

# 1. Rename the systems-report worksheet
$reportSheet = $wb.Worksheets.Item("تقرير الأنظمة")
$reportSheet.Name = "تقرير أنظمة البناء والتنظيم"

# 2-4. Update Assumptions sheet cells
$assump = $wb.Worksheets.Item("Assumptions")
$assump.Range("D3").Value = "Fund Name"
$assump.Range("E3").Value = "صندوق استثمار عقاري خاص"
$assump.Range("D4").Value = "Fund Type"

$wb.Save()
